$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.364.00"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "3.244.00"
$ws.Range("E3").Value = "  +2.37%  "
$ws.Range("D5").Value = "'605.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "'156.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.61%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.243.31"
$ws.Range("E8").Value = "  +2.42%  "
$ws.Range("D9").Value = "'0.546"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").Value = "'0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.81%  "
$ws.Range("D11").Value = "'5.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.04%  "
$ws.Range("D12").Value = "'0.497"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.67%  "
$ws.Range("D13").Value = "'0.0000269"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").Value = "'38.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").Value = "3.785.77"
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("D16").Value = "66.475.38"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "3.252.42"
$ws.Range("E18").Value = "  +3.00%  "
$ws.Range("D19").Value = "'0.114"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("D20").Value = "'504.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").Value = "'15.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "'0.749"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.84%  "
$ws.Range("D23").Value = "'8.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("D24").Value = "'14.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("D25").Value = "'87.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.76%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "'3.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").Value = "'9.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("D29").Value = "'2.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("E30").Value = "  +46.50%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'6.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.78%  "
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").Value = "'2.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.67%  "
$ws.Range("D33").Value = "'27.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("E35").Value = "  -3.42%  "
$ws.Range("D36").Value = "'6.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("D37").Value = "'3.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +18.65%  "
$ws.Range("D38").Value = "'55.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("D39").Value = "0.0₃0787"
$ws.Range("E39").Value = "  +14.80%  "
$ws.Range("D40").Value = "'493.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.57%  "
$ws.Range("D41").Value = "'0.0421"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("D43").Value = "'8.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.60%  "
$ws.Range("D44").Value = "'0.291"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.78%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.976.52"
$ws.Range("E45").Value = "  +5.41%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.17%  "
$ws.Range("D47").Value = "'28.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("D48").Value = "'2.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.69%  "
$ws.Range("E49").Value = "  +2.19%  "
$ws.Range("D51").Value = "'121.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.62%  "
